# Add a new "is_deleted" column (J) to the "All Products" sheet.
# Every existing product row gets a default value of 0 (not deleted),
# matching the new backend soft-delete flag used when a product is
# "deleted" (delisted instead of physically removed).

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("All Products")

# Determine the last used row reliably (row 55 in this workbook)
$lastRow = $ws.UsedRange.Rows.Count

# Header
$ws.Cells.Item(1, 10).Value = "is_deleted"

# Data rows: all existing product rows default to not-deleted (0)
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 10).Value = 0
}

$ws.Range("J2").Select()
